$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Output $ws.Name
